$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.376.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.870.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06462'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.97'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07788'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.15'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.865.99'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7240'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.139'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '281.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.367.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.0000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007499'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.114.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.253'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.240'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09612'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.485'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.230'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.113'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04816'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6898'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.713'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01890'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.820'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.930'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8273'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.649'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.959'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '898.61'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.36%  '
$ws.Range('E51').Value = '  +0.51%  '

Write-Output "Applied changes successfully"